$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 326.33334
$ws.Range("I2").Value = 151.77777
$ws.Range("J2").Value = 850
$ws.Range("K2").Value = 151.77777
$ws.Range("L2").Value = 850
$ws.Range("M2").Value = -38.77777
$ws.Range("N2").Value = -1076

# Row 70
$ws.Range("H70").Value = 1569.1428
$ws.Range("I70").Value = 896
$ws.Range("J70").Value = 1905.7142
$ws.Range("K70").Value = 2688
$ws.Range("L70").Value = 5717.142599999999
$ws.Range("M70").Value = -2418
$ws.Range("N70").Value = -6257.142599999999

# Row 73
$ws.Range("H73").Value = 1569.1428
$ws.Range("I73").Value = 896
$ws.Range("J73").Value = 1905.7142
$ws.Range("K73").Value = 2688
$ws.Range("L73").Value = 5717.142599999999
$ws.Range("M73").Value = -1752
$ws.Range("N73").Value = -7589.142599999999

# Row 80
$ws.Range("H80").Value = 640.1739
$ws.Range("I80").Value = 623.55554
$ws.Range("J80").Value = 700
$ws.Range("K80").Value = 1870.66662
$ws.Range("L80").Value = 2100
$ws.Range("M80").Value = -872.66662
$ws.Range("N80").Value = -4096

# Row 83
$ws.Range("H83").Value = 640.1739
$ws.Range("I83").Value = 623.55554
$ws.Range("J83").Value = 700
$ws.Range("K83").Value = 5611.99986
$ws.Range("L83").Value = 6300
$ws.Range("M83").Value = -619.9998599999999
$ws.Range("N83").Value = -16284

# Row 95
$ws.Range("H95").Value = 35999.8
$ws.Range("J95").Value = 35999.8
$ws.Range("L95").Value = 35999.8
$ws.Range("N95").Value = -41491.8

# Row 111
$ws.Range("H111").Value = 861.46155
$ws.Range("I111").Value = 957
$ws.Range("J111").Value = 750
$ws.Range("K111").Value = 2871
$ws.Range("L111").Value = 2250
$ws.Range("M111").Value = 196
$ws.Range("N111").Value = -8384

# Row 116
$ws.Range("H116").Value = 4634.6523
$ws.Range("I116").Value = 3287.625
$ws.Range("J116").Value = 5353.067
$ws.Range("K116").Value = 3287.625
$ws.Range("L116").Value = 5353.067
$ws.Range("M116").Value = 154.375
$ws.Range("N116").Value = -12237.067

# Row 141
$ws.Range("H141").Value = 9925.637000000001
$ws.Range("I141").Value = 3350.6
$ws.Range("J141").Value = 24015
$ws.Range("K141").Value = 10051.8
$ws.Range("L141").Value = 72045
$ws.Range("M141").Value = -4871.799999999999
$ws.Range("N141").Value = -82405

$ws = $wb.Worksheets.Item("ARM")
# Row 95
$ws.Range("H95").Value = 21235.1
$ws.Range("J95").Value = 21235.1
$ws.Range("L95").Value = 21235.1
$ws.Range("N95").Value = -26727.1

# Row 122
$ws.Range("H122").Value = 2602.8333
$ws.Range("I122").Value = 1021.2727
$ws.Range("J122").Value = 20000
$ws.Range("K122").Value = 3063.8181
$ws.Range("L122").Value = 60000
$ws.Range("M122").Value = -613.8181
$ws.Range("N122").Value = -64900

# Row 132
$ws.Range("H132").Value = 8733.684999999999
$ws.Range("I132").Value = 3991
$ws.Range("J132").Value = 12182.909
$ws.Range("K132").Value = 11973
$ws.Range("L132").Value = 36548.727
$ws.Range("M132").Value = -9443
$ws.Range("N132").Value = -41608.727

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 45498.5
$ws.Range("J35").Value = 46798.2
$ws.Range("L35").Value = 46798.2
$ws.Range("N35").Value = -47418.2

# Row 76
$ws.Range("H76").Value = 40314
$ws.Range("J76").Value = 40314
$ws.Range("L76").Value = 40314
$ws.Range("N76").Value = -40944

# Row 79
$ws.Range("H79").Value = 40314
$ws.Range("J79").Value = 40314
$ws.Range("L79").Value = 40314
$ws.Range("N79").Value = -42498

# Row 92
$ws.Range("H92").Value = 20500
$ws.Range("J92").Value = 20500
$ws.Range("L92").Value = 20500
$ws.Range("N92").Value = -25492

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 3290.9614
$ws.Range("I132").Value = 2692.8572
$ws.Range("J132").Value = 3511.3157
$ws.Range("K132").Value = 8078.571599999999
$ws.Range("L132").Value = 10533.9471
$ws.Range("M132").Value = -5548.571599999999
$ws.Range("N132").Value = -15593.9471

$ws = $wb.Worksheets.Item("GSM")
# Row 95
$ws.Range("H95").Value = 15747
$ws.Range("J95").Value = 15747
$ws.Range("L95").Value = 15747
$ws.Range("N95").Value = -21239

# Row 102
$ws.Range("H102").Value = 3971847.2
$ws.Range("I102").Value = 10207280
$ws.Range("K102").Value = 10207280
$ws.Range("M102").Value = -10205658

# Row 122
$ws.Range("H122").Value = 5002500
$ws.Range("I122").Value = 5002500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 15007500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -15005050
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 50004524
$ws.Range("I40").Value = 71431670
$ws.Range("J40").Value = 7831.6665
$ws.Range("K40").Value = 71431670
$ws.Range("L40").Value = 7831.6665
$ws.Range("M40").Value = -71431534
$ws.Range("N40").Value = -8103.6665

# Row 122
$ws.Range("H122").Value = 4875.3335
$ws.Range("I122").Value = 5181.6
$ws.Range("J122").Value = 4492.5
$ws.Range("K122").Value = 15544.8
$ws.Range("L122").Value = 13477.5
$ws.Range("M122").Value = -13094.8
$ws.Range("N122").Value = -18377.5

$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 27500
$ws.Range("J16").Value = 27500
$ws.Range("L16").Value = 27500
$ws.Range("N16").Value = -28084

# Row 81
$ws.Range("H81").Value = 2178.111
$ws.Range("I81").Value = 3150.375
$ws.Range("J81").Value = 1400.3
$ws.Range("K81").Value = 6300.75
$ws.Range("L81").Value = 2800.6
$ws.Range("M81").Value = -5239.75
$ws.Range("N81").Value = -4922.6

# Row 84
$ws.Range("H84").Value = 2178.111
$ws.Range("I84").Value = 3150.375
$ws.Range("J84").Value = 1400.3
$ws.Range("K84").Value = 31503.75
$ws.Range("L84").Value = 14003
$ws.Range("M84").Value = -26199.75
$ws.Range("N84").Value = -24611

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 136
$ws.Range("H136").Value = 4718611.5
$ws.Range("I136").Value = 8065651.5
$ws.Range("J136").Value = 2327.7273
$ws.Range("K136").Value = 24196954.5
$ws.Range("L136").Value = 6983.1819
$ws.Range("M136").Value = -24194404.5
$ws.Range("N136").Value = -12083.1819
